# Reference Genome Assembly Test suite
# Rename the "Cases" tab entry to "Participants" on the startup sheet,
# and refresh the view/formatting the way Excel does when the sheet is
# next edited and saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Core content change: A2 "CasesTab" -> "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Row heights grow to match the new default row height / font metrics
# used when the workbook was next recalculated (wrapped text rows).
$ws.Rows("2:2").RowHeight = 165
$ws.Rows("3:3").RowHeight = 180
$ws.Rows("4:4").RowHeight = 210

# Column A widens (best-fit) to accommodate the longer "ParticipantsTab" text.
$ws.Columns("A:A").ColumnWidth = 15.73

# Selection moves from a single cell to the whole sheet (select-all).
$ws.Cells.Select() | Out-Null
